$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set Execute column (B) to "Yes" for all test rows (3-25) except row 20 which was
# already "Yes". Row 2 stays "No" (unchanged by this edit).
$rowsToEnable = 3..19 + 21..25
foreach ($r in $rowsToEnable) {
    $ws.Cells.Item($r, 2).Value = "Yes"
}

# Extend the Keywords2 test parameters for "Verify Merchant Activity Links" (row 21)
$ws.Cells.Item(21, 9).Value = "coyni.merchant.tests.MerchantActivityTest,`ntestMerchantActivityLinks,`n-pheading,`n-ptransactionHeading,`n-ppayOutHistoryHeading,`n-preserveHistoryHeading"

# Update the view selection to reflect the final saved cursor state
$ws.Range("B5").Select()
